$wb = $excel.ActiveWorkbook

# 1. Update the comment on O1 (single_file_export_format) in "Export as TSV" sheet
$wsMain = $wb.Worksheets.Item("Export as TSV")
$comment = $wsMain.Range("O1").Comment
$comment.Text("The format in which each single imaging file will be exported. (Example: DICOM, tiff, avi, etc.)")

# 2. Update the "quality_view list" sheet: replace "high" with three new options
$wsQuality = $wb.Worksheets.Item("quality_view list")
$wsQuality.Range("A1").Value = "high confidence/optimal"
$wsQuality.Range("A2").Value = "low confidence/sub-optimal"
$wsQuality.Range("A3").Value = "no confidence"

# 3. Update the data validation on column AE (quality_view) to reference the
#    expanded list range and updated error message
$dv = $wsMain.Range("AE2:AE1048576").Validation
$dv.Formula1 = "'quality_view list'!`$A`$1:`$A`$3"
$dv.ErrorMessage = "Value must be one of: high confidence/optimal / low confidence/sub-optimal / no confidence."
